$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customers")

# Duplicate the header row plus the first 3 data rows into a brand-new sheet
# placed right after "Customers" (Excel auto-names it "Customers1").
$new = $wb.Worksheets.Add($null, $ws)
$new.Name = "Customers1"
$ws.Range("A1:G4").Copy($new.Range("A1:G4"))

# Freeze panes the same way as the source sheet (6 cols / 1 row frozen).
$new.Activate()
$new.Range("G2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Leave the pasted block selected, matching the paste-then-select state.
$new.Range("A2:G4").Select()
